# Auto-generated Excel COM-interop script
# Applies a scheduled market-data refresh to the Leve profit sheets
# (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*)

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 42
$ws.Range("H42").Value = 545.5
$ws.Range("I42").Value = 62.57143
$ws.Range("J42").Value = 1221.6
$ws.Range("K42").Value = 187.71429
$ws.Range("L42").Value = 3664.8
$ws.Range("M42").Value = 42.28570999999999
$ws.Range("N42").Value = -4124.799999999999
# Row 112
$ws.Range("H112").Value = 1385.2084
$ws.Range("J112").Value = 1562.25
$ws.Range("L112").Value = 4686.75
$ws.Range("N112").Value = -6902.75
# Row 125
$ws.Range("H125").Value = 622.9091
$ws.Range("I125").Value = 550
$ws.Range("J125").Value = 639.1111
$ws.Range("K125").Value = 4950
$ws.Range("L125").Value = 5751.9999
$ws.Range("M125").Value = -2490
$ws.Range("N125").Value = -10671.9999
# Row 138
$ws.Range("H138").Value = 1863.7377
$ws.Range("I138").Value = 653.0857
$ws.Range("J138").Value = 3493.4614
$ws.Range("K138").Value = 1959.2571
$ws.Range("L138").Value = 10480.3842
$ws.Range("M138").Value = 3180.7429
$ws.Range("N138").Value = -20760.3842

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1111.4333
$ws.Range("I61").Value = 1002.5263
$ws.Range("J61").Value = 1299.5454
$ws.Range("K61").Value = 1002.5263
$ws.Range("L61").Value = 1299.5454
$ws.Range("M61").Value = -790.5263
$ws.Range("N61").Value = -1723.5454
# Row 123
$ws.Range("H123").Value = 45958.57
$ws.Range("J123").Value = 45958.57
$ws.Range("L123").Value = 45958.57
$ws.Range("N123").Value = -55758.57
# Row 136
$ws.Range("H136").Value = 1111.4333
$ws.Range("I136").Value = 1002.5263
$ws.Range("J136").Value = 1299.5454
$ws.Range("K136").Value = 3007.5789
$ws.Range("L136").Value = 3898.6362
$ws.Range("M136").Value = -457.5789
$ws.Range("N136").Value = -8998.636200000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2058.6538
$ws.Range("I134").Value = 1773.2273
$ws.Range("J134").Value = 3628.5
$ws.Range("K134").Value = 5319.6819
$ws.Range("L134").Value = 10885.5
$ws.Range("M134").Value = -2784.6819
$ws.Range("N134").Value = -15955.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4188957.8
$ws.Range("I31").Value = 7179499.5
$ws.Range("K31").Value = 7179499.5
$ws.Range("M31").Value = -7179204.5
# Row 34
$ws.Range("H34").Value = 4188957.8
$ws.Range("I34").Value = 7179499.5
$ws.Range("K34").Value = 7179499.5
$ws.Range("M34").Value = -7179297.5
# Row 58
$ws.Range("H58").Value = 1614.174
$ws.Range("I58").Value = 877
$ws.Range("J58").Value = 2289.9167
$ws.Range("K58").Value = 877
$ws.Range("L58").Value = 2289.9167
$ws.Range("M58").Value = -674
$ws.Range("N58").Value = -2695.9167
# Row 134
$ws.Range("H134").Value = 3164.6155
$ws.Range("I134").Value = 3367.6
$ws.Range("J134").Value = 2488
$ws.Range("K134").Value = 10102.8
$ws.Range("L134").Value = 7464
$ws.Range("M134").Value = -7567.799999999999
$ws.Range("N134").Value = -12534
# Row 136
$ws.Range("H136").Value = 1614.174
$ws.Range("I136").Value = 877
$ws.Range("J136").Value = 2289.9167
$ws.Range("K136").Value = 2631
$ws.Range("L136").Value = 6869.750100000001
$ws.Range("M136").Value = -81
$ws.Range("N136").Value = -11969.7501

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 127
$ws.Range("H127").Value = 1052.3334
$ws.Range("J127").Value = 1052.3334
$ws.Range("L127").Value = 3157.0002
$ws.Range("N127").Value = -13077.0002
# Row 131
$ws.Range("H131").Value = 1170442.8
$ws.Range("I131").Value = 3509035.5
$ws.Range("J131").Value = 1146.3422
$ws.Range("K131").Value = 10527106.5
$ws.Range("L131").Value = 3439.0266
$ws.Range("M131").Value = -10522066.5
$ws.Range("N131").Value = -13519.0266

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2203.2122
$ws.Range("I80").Value = 2152.6667
$ws.Range("J80").Value = 2291.6667
$ws.Range("K80").Value = 2152.6667
$ws.Range("L80").Value = 2291.6667
$ws.Range("M80").Value = -1154.6667
$ws.Range("N80").Value = -4287.6667
# Row 83
$ws.Range("H83").Value = 2203.2122
$ws.Range("I83").Value = 2152.6667
$ws.Range("J83").Value = 2291.6667
$ws.Range("K83").Value = 10763.3335
$ws.Range("L83").Value = 11458.3335
$ws.Range("M83").Value = -5771.333500000001
$ws.Range("N83").Value = -21442.3335
# Row 123
$ws.Range("H123").Value = 17078.625
$ws.Range("J123").Value = 17078.625
$ws.Range("L123").Value = 17078.625
$ws.Range("N123").Value = -21978.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 55556544
$ws.Range("I93").Value = 62500988
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 62500988
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -62499740
$ws.Range("N93").Value = -3496
# Row 100
$ws.Range("H100").Value = 1565.5555
$ws.Range("I100").Value = 1516.3636
$ws.Range("J100").Value = 1642.8572
$ws.Range("K100").Value = 1516.3636
$ws.Range("L100").Value = 1642.8572
$ws.Range("M100").Value = -975.3635999999999
$ws.Range("N100").Value = -2724.8572
# Row 123
$ws.Range("H123").Value = 17593.545
$ws.Range("J123").Value = 17593.545
$ws.Range("L123").Value = 17593.545
$ws.Range("N123").Value = -27393.545
# Row 124
$ws.Range("H124").Value = 30476.334
$ws.Range("J124").Value = 30476.334
$ws.Range("L124").Value = 30476.334
$ws.Range("N124").Value = -40296.334
# Row 125
$ws.Range("H125").Value = 40838.332
$ws.Range("J125").Value = 40838.332
$ws.Range("L125").Value = 40838.332
$ws.Range("N125").Value = -50678.332
# Row 127
$ws.Range("H127").Value = 29505
$ws.Range("J127").Value = 29505
$ws.Range("L127").Value = 29505
$ws.Range("N127").Value = -39425
# Row 128
$ws.Range("H128").Value = 22222
$ws.Range("J128").Value = 22222
$ws.Range("L128").Value = 22222
$ws.Range("N128").Value = -32182
# Row 129
$ws.Range("H129").Value = 38000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 38000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 38000
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -48000
# Row 130
$ws.Range("H130").Value = 25000
$ws.Range("J130").Value = 25000
$ws.Range("L130").Value = 25000
$ws.Range("N130").Value = -35040
# Row 131
$ws.Range("H131").Value = 38775
$ws.Range("J131").Value = 38775
$ws.Range("L131").Value = 38775
$ws.Range("N131").Value = -48855

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2284
$ws.Range("I122").Value = 1901.3334
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 5704.0002
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").Value = -3254.0002
$ws.Range("N122").Value = -12900.0001
# Row 123
$ws.Range("H123").Value = 46238.395
$ws.Range("J123").Value = 46238.395
$ws.Range("L123").Value = 46238.395
$ws.Range("N123").Value = -56038.395
# Row 136
$ws.Range("H136").Value = 857.2033699999999
$ws.Range("I136").Value = 453.3654
$ws.Range("K136").Value = 1360.0962
$ws.Range("M136").Value = 1189.9038
